$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginData")
$ws.Activate()

# Update the runMode for row 3 from "Y" to "N"
$ws.Range("D3").Value = "N"

# Update the password/data value for row 4 (leading apostrophe keeps the
# existing "quote prefix" text style on the cell, same as before the edit)
$ws.Range("C4").Value = "'f1234567890k"

# Move the selection to B9 to match the saved cursor position
$ws.Range("B9").Select()
